$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Summer
$ws.Range("G6").Value = "Summer"

# Rows 7-12: Spring
$ws.Range("G7").Value = "Spring"
$ws.Range("G8").Value = "Spring"
$ws.Range("G9").Value = "Spring"
$ws.Range("G10").Value = "Spring"
$ws.Range("G11").Value = "Spring"
$ws.Range("G12").Value = "Spring"

# Rows 13-18: Fall
$ws.Range("G13").Value = "Fall"
$ws.Range("G14").Value = "Fall"
$ws.Range("G15").Value = "Fall"
$ws.Range("G16").Value = "Fall"
$ws.Range("G17").Value = "Fall"
$ws.Range("G18").Value = "Fall"

# Update selection to match the saved cursor position
$ws.Range("A6").Select() | Out-Null
